# Auto-generated edit script: updates crypto price/volume cells per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (losing the exact original formatting, e.g. trailing zeros) are written
# via a temporary "@" (Text) number format, then ClearFormats() restores the
# original (unstyled) cell style so no stray formatting diff is introduced.

$ws.Range("D2").Value = '30.344.20'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '1.869.89'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.72'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4700'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.25%  '
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06444'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.04'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07756'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").Value = '1.869.97'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.89'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7212'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.123'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.24'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("D17").Value = '30.326.02'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.96'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007476'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").Value = '2.112.65'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.213'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.25'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.031'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.64'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.869'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.12%  '
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09590'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.65%  '
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.192'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.090'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04802'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6879'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.716'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01876'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("E39").Value = '  +1.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.196'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.20'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4214'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.928'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9992'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8278'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.53'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.534'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.23'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.933'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '900.88'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05718'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.90%  '

Write-Output "Updated 91 cells (35 via text-format workaround, 56 direct)"
